# Swap the data in columns A, B, E, F, G, H, Q, R between row 2 and row 3.
# These columns identify the taxon/observation details that got reordered;
# all other columns share identical values between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
